$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.504.71"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "3.942.20"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "492.14"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.78"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  +4.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000349"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.08"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "4.569.15"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "3.938.56"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  -4.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.90"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "69.454.22"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.58"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.55"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.49"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.05"
$ws.Range("E25").Value = "  +9.48%  "
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.13"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.34"
$ws.Range("E28").Value = "  -4.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.67"
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "707.79"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.474"
$ws.Range("E34").Value = "  +26.10%  "
$ws.Range("D35").Value = "0.0₃0909"
$ws.Range("E35").Value = "  -4.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "61.70"
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").Value = "  +4.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "40.79"
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0491"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.94"
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.08"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "0.0₆0367"
$ws.Range("E47").Value = "  +11.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.34"
$ws.Range("E48").Value = "  +6.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.08"
$ws.Range("E49").Value = "  +8.35%  "
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.14"
$ws.Range("E51").Value = "  -1.49%  "
